$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11th column). This shifts everything
# from K onward one column to the right, and Excel automatically grows any
# merged ranges / defined ranges that span the insertion point.
$ws.Columns("K:K").Insert()

Write-Host "done"
